$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F - "Trening" (training split) header, styled like the other headers
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Trening"

# Convert the Timestamp column (A2:A7) from text strings into real Excel
# date serial values, and mark each row with its training part ("Gra").
$serials = @(45687.52221782407, 45687.52894004629, 45687.52991226852, 45687.52221666666, 45687.52893888889, 45687.5299087963)

for ($i = 0; $i -lt $serials.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $serials[$i]
    $ws.Cells.Item($row, 6).Value = "Gra"
}

# Apply the custom date/time display format to the timestamp column.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
